# initial version of stimuli order
# Rewrites the per-trial stimulus assignment (B: index into the 0..123
# shuffled pool, C: image, D: word, E: category) for rows 2-33 to the
# new randomized ordering.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, B (number), C (image), D (word), E (category)
$rows = @(
    @(2,   122, 'face/face100.png', 'nehmen',   'face'),
    @(3,   108, 'face/face120.png', 'fliehen',  'face'),
    @(4,    14, 'face/face079.png', 'füttern',  'face'),
    @(5,    63, 'car/car074.png',   'raten',    'car'),
    @(6,    41, 'car/car099.png',   'loben',    'car'),
    @(7,   107, 'car/car122.png',   'husten',   'car'),
    @(8,    70, 'car/car107.png',   'gelten',   'car'),
    @(9,    36, 'car/car081.png',   'opfern',   'car'),
    @(10,  103, 'face/face072.png', 'regnen',   'face'),
    @(11,   96, 'car/car101.png',   'pflegen',  'car'),
    @(12,   99, 'car/car120.png',   'schätzen', 'car'),
    @(13,   92, 'car/car083.png',   'wiegen',   'car'),
    @(14,   71, 'face/face104.png', 'dauern',   'face'),
    @(15,   64, 'car/car069.png',   'tagen',    'car'),
    @(16,  117, 'face/face105.png', 'enden',    'face'),
    @(17,  114, 'car/car075.png',   'starten',  'car'),
    @(18,   10, 'face/face091.png', 'fühlen',   'face'),
    @(19,  101, 'car/car064.png',   'biegen',   'car'),
    @(20,  112, 'face/face123.png', 'liefern',  'face'),
    @(21,   69, 'face/face109.png', 'währen',   'face'),
    @(22,   33, 'car/car078.png',   'hupen',    'car'),
    @(23,   27, 'car/car072.png',   'antun',    'car'),
    @(24,   26, 'face/face064.png', 'mieten',   'face'),
    @(25,   86, 'face/face098.png', 'kaufen',   'face'),
    @(26,   39, 'face/face107.png', 'stechen',  'face'),
    @(27,  118, 'car/car116.png',   'ehren',    'car'),
    @(28,   28, 'face/face096.png', 'töten',    'face'),
    @(29,   19, 'face/face083.png', 'wenden',   'face'),
    @(30,   72, 'car/car110.png',   'backen',   'car'),
    @(31,   37, 'face/face102.png', 'sieben',   'face'),
    @(32,   30, 'face/face090.png', 'rasen',    'face'),
    @(33,   67, 'car/car065.png',   'schenken', 'car')
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 2).Value = $r[1]
    $ws.Cells.Item($rowNum, 3).Value = $r[2]
    $ws.Cells.Item($rowNum, 4).Value = $r[3]
    $ws.Cells.Item($rowNum, 5).Value = $r[4]
}
